$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C width
$ws.Range("C1").ColumnWidth = 23.8

# E2 wage change
$ws.Range("E2").Value = 7

# Apply body-row format (style 5/6) to A3:F22 by copying from an existing body row
$ws.Range("A12:F12").Copy()
$ws.Range("A3:F22").PasteSpecial(-4122)

# Row heights
$ws.Rows.Item(2).RowHeight = 20.25
$ws.Rows.Item(3).RowHeight = 20.05
$ws.Rows.Item(4).RowHeight = 20.05
$ws.Rows.Item(5).RowHeight = 20.05
$ws.Rows.Item(6).RowHeight = 20.05
$ws.Rows.Item(7).RowHeight = 20.05
$ws.Rows.Item(8).RowHeight = 32.05
$ws.Rows.Item(9).RowHeight = 32.05
$ws.Rows.Item(10).RowHeight = 20.05
$ws.Rows.Item(11).RowHeight = 20.05
$ws.Rows.Item(12).RowHeight = 20.05
$ws.Rows.Item(13).RowHeight = 20.05
$ws.Rows.Item(14).RowHeight = 20.05
$ws.Rows.Item(15).RowHeight = 20.05
$ws.Rows.Item(16).RowHeight = 20.05
$ws.Rows.Item(17).RowHeight = 20.05
$ws.Rows.Item(18).RowHeight = 20.05
$ws.Rows.Item(19).RowHeight = 20.05
$ws.Rows.Item(20).RowHeight = 20.05
$ws.Rows.Item(21).RowHeight = 20.05
$ws.Rows.Item(22).RowHeight = 20.05

# Cell values
$ws.Range("A2").Value = "Finance"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "Admin Assistant"
$ws.Range("D2").Value = 12
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = "Monday, Tuesday, Wednesday, Thursday, Friday"
$ws.Range("A3").Value = "Finance"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "Invoice Processor"
$ws.Range("D3").Value = 10.5
$ws.Range("E3").Value = 7.5
$ws.Range("F3").Value = "Monday, Tuesday, Wednesday, Thursday, Friday"
$ws.Range("A4").Value = "Finance"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = "Expenses Assistant"
$ws.Range("D4").Value = 15
$ws.Range("E4").Value = 7.5
$ws.Range("F4").Value = "Monday, Tuesday, Wednesday, Thursday, Friday"
$ws.Range("A5").Value = "Finance"
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = "Payments Assistant"
$ws.Range("D5").Value = 17
$ws.Range("E5").Value = 7.5
$ws.Range("F5").Value = "Monday, Tuesday, Wednesday, Thursday, Friday"
$ws.Range("A6").Value = "Finance"
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = "Accounts Receivable Assistant"
$ws.Range("D6").Value = 17
$ws.Range("E6").Value = 7.5
$ws.Range("F6").Value = "Monday, Tuesday, Wednesday, Thursday, Friday"
$ws.Range("A7").Value = "Finance"
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = "Accounts Payable Assistant"
$ws.Range("D7").Value = 17
$ws.Range("E7").Value = 7.5
$ws.Range("F7").Value = "Monday, Tuesday, Wednesday, Thursday, Friday"
$ws.Range("A8").Value = "Finance"
$ws.Range("B8").Value = 4
$ws.Range("C8").Value = "Senior Accounts Payable Assistant"
$ws.Range("D8").Value = 20
$ws.Range("E8").Value = 7.5
$ws.Range("F8").Value = "Monday, Tuesday, Wednesday, Thursday, Friday"
$ws.Range("A9").Value = "Finance"
$ws.Range("B9").Value = 4
$ws.Range("C9").Value = "Senior Accounts Receivable Assistant"
$ws.Range("D9").Value = 20
$ws.Range("E9").Value = 7.5
$ws.Range("F9").Value = "Monday, Tuesday, Wednesday, Thursday, Friday"
$ws.Range("A10").Value = "Finance"
$ws.Range("B10").Value = 4
$ws.Range("C10").Value = "Accounts Assistant"
$ws.Range("D10").Value = 22
$ws.Range("E10").Value = 7.5
$ws.Range("F10").Value = "Monday, Tuesday, Wednesday, Thursday, Friday"
$ws.Range("A11").Value = "Finance"
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = "Part Qualified Accountant"
$ws.Range("D11").Value = 32
$ws.Range("E11").Value = 7.5
$ws.Range("F11").Value = "Monday, Tuesday, Wednesday, Thursday, Friday"
$ws.Range("A12").Value = "Finance"
$ws.Range("B12").Value = 6
$ws.Range("C12").Value = "Chartered Accountant"
$ws.Range("D12").Value = 52
$ws.Range("E12").Value = 7
$ws.Range("F12").Value = "Monday, Tuesday, Wednesday, Thursday, Friday"
$ws.Range("A13").Value = "Finance"
$ws.Range("B13").Value = 7
$ws.Range("C13").Value = "Finance Manager"
$ws.Range("D13").Value = 68
$ws.Range("E13").Value = 7
$ws.Range("F13").Value = "Monday, Tuesday, Wednesday, Thursday, Friday"
$ws.Range("A14").Value = "Finance"
$ws.Range("B14").Value = 8
$ws.Range("C14").Value = "Financial Controller"
$ws.Range("D14").Value = 80
$ws.Range("E14").Value = 7
$ws.Range("F14").Value = "Monday, Tuesday, Wednesday, Thursday, Friday"
$ws.Range("A15").Value = "Finance"
$ws.Range("B15").Value = 9
$ws.Range("C15").Value = "CFO"
$ws.Range("D15").Value = 95
$ws.Range("E15").Value = 6.5
$ws.Range("F15").Value = "Monday, Tuesday, Wednesday, Thursday, Friday"
$ws.Range("A16").Value = "Finance"
$ws.Range("B16").Value = 10
$ws.Range("C16").Value = "CEO"
$ws.Range("D16").Value = 110
$ws.Range("E16").Value = 6.5
$ws.Range("F16").Value = "Monday, Tuesday, Wednesday, Thursday"
$ws.Range("A17").Value = "Culinary"
$ws.Range("B17").Value = 1
$ws.Range("C17").Value = "Dishwasher"
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = 8
$ws.Range("F17").Value = "Tuesday, Wednesday, Friday, Saturday, Sunday"
$ws.Range("A18").Value = "Culinary"
$ws.Range("B18").Value = 2
$ws.Range("C18").Value = "Salad Prep"
$ws.Range("D18").Value = 12
$ws.Range("E18").Value = 8
$ws.Range("F18").Value = "Tuesday, Wednesday, Friday, Saturday, Sunday"
$ws.Range("A19").Value = "Culinary"
$ws.Range("B19").Value = 3
$ws.Range("C19").Value = "Line Cook"
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = 8
$ws.Range("F19").Value = "Tuesday, Wednesday, Friday, Saturday, Sunday"
$ws.Range("A20").Value = "Culinary"
$ws.Range("B20").Value = 4
$ws.Range("C20").Value = "Dessert Specialist"
$ws.Range("D20").Value = 18
$ws.Range("E20").Value = 8
$ws.Range("F20").Value = "Tuesday, Wednesday, Friday, Saturday, Sunday"
$ws.Range("A21").Value = "Culinary"
$ws.Range("B21").Value = 5
$ws.Range("C21").Value = "Sous Chef"
$ws.Range("D21").Value = 28
$ws.Range("E21").Value = 8
$ws.Range("F21").Value = "Tuesday, Wednesday, Friday, Saturday, Sunday"
$ws.Range("A22").Value = "Culinary"
$ws.Range("B22").Value = 6
$ws.Range("C22").Value = "Head Chef"
$ws.Range("D22").Value = 40
$ws.Range("E22").Value = 8
$ws.Range("F22").Value = "Tuesday, Wednesday, Friday, Saturday, Sunday"
